$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Select the cell (mirrors the user clicking on E8 before editing it)
$ws.Range("E8").Select()

# Update the greeting text from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"
